$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Junio de 2020 a las 02:16"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 2142017
$ws.Range("C4").Value = 25095
$ws.Range("D4").Value = 851135
$ws.Range("E4").Value = 1173356
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 701
$ws.Range("H4").Value = 117526

# Row 5: Brasil
$ws.Range("A5").Value = "Brasil"
$ws.Range("B5").Value = 850796
$ws.Range("C5").Value = 20894
$ws.Range("D5").Value = 427610
$ws.Range("E5").Value = 380395
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 890
$ws.Range("H5").Value = 42791

# Row 20: Canada
$ws.Range("A20").Value = "Canada"
$ws.Range("B20").Value = 98410
$ws.Range("C20").Value = 467
$ws.Range("D20").Value = 59354
$ws.Range("E20").Value = 30949
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 58
$ws.Range("H20").Value = 8107

# Row 33: Singapur
$ws.Range("A33").Value = "Singapur"
$ws.Range("B33").Value = 40197
$ws.Range("C33").Value = 347
$ws.Range("D33").Value = 28808
$ws.Range("E33").Value = 11363
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 26

# Row 47: Panama
$ws.Range("A47").Value = "Panama"
$ws.Range("B47").Value = 20059
$ws.Range("C47").Value = 848
$ws.Range("D47").Value = 13759
$ws.Range("E47").Value = 5871
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 429

# Row 93: Venezuela
$ws.Range("A93").Value = "Venezuela"
$ws.Range("B93").Value = 2904
$ws.Range("C93").Value = 25
$ws.Range("D93").Value = 487
$ws.Range("E93").Value = 2393
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 24

# Row 94: Bosnia y Herzegovina
$ws.Range("A94").Value = "Bosnia y Herzegovina"
$ws.Range("B94").Value = 2893
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 2119
$ws.Range("E94").Value = 611
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 163

# Row 148: Estado de Palestina
$ws.Range("A148").Value = "Estado de Palestina"
$ws.Range("B148").Value = 489
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 415
$ws.Range("E148").Value = 71
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 3

# Row 153: Libia
$ws.Range("A153").Value = "Libia"
$ws.Range("B153").Value = 418
$ws.Range("C153").Value = 9
$ws.Range("D153").Value = 62
$ws.Range("E153").Value = 348
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 2
$ws.Range("H153").Value = 8

# Row 154: Benin
$ws.Range("A154").Value = "Benin"
$ws.Range("B154").Value = 412
$ws.Range("C154").Value = 24
$ws.Range("D154").Value = 222
$ws.Range("E154").Value = 184
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 6

# Row 183: Eritrea
$ws.Range("A183").Value = "Eritrea"
$ws.Range("B183").Value = 65
$ws.Range("C183").Value = 24
$ws.Range("D183").Value = 39
$ws.Range("E183").Value = 26
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

# Row 184: Butan
$ws.Range("A184").Value = "Butan"
$ws.Range("B184").Value = 62
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 20
$ws.Range("E184").Value = 42
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

# Row 185: Botsuana
$ws.Range("A185").Value = "Botsuana"
$ws.Range("B185").Value = 60
$ws.Range("C185").Value = 12
$ws.Range("D185").Value = 24
$ws.Range("E185").Value = 35
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 1

# Row 186: Polinesia Francesa
$ws.Range("A186").Value = "Polinesia Francesa"
$ws.Range("B186").Value = 60
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 60
$ws.Range("E186").Value = 0
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

# Row 187: Macao
$ws.Range("A187").Value = "Macao"
$ws.Range("B187").Value = 45
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 45
$ws.Range("E187").Value = 0
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0

# Row 188: San Martin (Parte Francesa)
$ws.Range("A188").Value = "San Martin (Parte Francesa)"
$ws.Range("B188").Value = 42
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 36
$ws.Range("E188").Value = 3
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 3

# Row 210: Seychelles
$ws.Range("A210").Value = "Seychelles"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 11
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211: Montserrat
$ws.Range("A211").Value = "Montserrat"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 10
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 1
